$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 75: B value changes
$ws.Range("B75").Value = 411000000

# Row 76: A and B change
$ws.Range("A76").Value = "aunties"
$ws.Range("B76").Value = 54400000

# Row 77: A and B change
$ws.Range("A77").Value = "cheer"
$ws.Range("B77").Value = 186000000

# Row 78: A and B change
$ws.Range("A78").Value = "variables"
$ws.Range("B78").Value = 767000000

# Row 79: B value changes
$ws.Range("B79").Value = 411000000

# Row 80: A and B change
$ws.Range("A80").Value = "variables"
$ws.Range("B80").Value = 767000000

# Row 81: A and B change
$ws.Range("A81").Value = "question"
$ws.Range("B81").Value = 2510000000

# Row 82: A and B change
$ws.Range("A82").Value = "ideologies"
$ws.Range("B82").Value = 21700000

# Row 83: B value changes
$ws.Range("B83").Value = 404000000

# Row 84: A and B change
$ws.Range("A84").Value = "excite"
$ws.Range("B84").Value = 103000000

# Row 85: A and B change
$ws.Range("A85").Value = "birthday"
$ws.Range("B85").Value = 1810000000

# Row 86: A and B change
$ws.Range("A86").Value = "different"
$ws.Range("B86").Value = 4580000000

# New rows 87-91
$ws.Range("A87").Value = "tesla"
$ws.Range("B87").Value = 388000000

$ws.Range("A88").Value = "excite"
$ws.Range("B88").Value = 103000000

$ws.Range("A89").Value = "clarity"
$ws.Range("B89").Value = 205000000

$ws.Range("A90").Value = "wish"
$ws.Range("B90").Value = 2590000000

$ws.Range("A91").Value = "tesla"
$ws.Range("B91").Value = 388000000
